# Actualización 11 de Mayo - Mañana
# "Rescatables" sheet: two new students are added (rescue/remedial exam
# candidates for group 2BLCM) ahead of the existing entry, and the
# existing student's Reprobadas count is corrected from 2 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# New row 2 - ANTONIO GARCIA ISRAEL (2BLCM)
$ws.Cells.Item(2,1).Value = 20330051920284
$ws.Cells.Item(2,2).Value = "ANTONIO"
$ws.Cells.Item(2,3).Value = "GARCIA"
$ws.Cells.Item(2,4).Value = "ISRAEL"
$ws.Cells.Item(2,5).Value = "TOMA MUESTRAS BIOLÓGICAS"
$ws.Cells.Item(2,6).Value = "2BLCM"
$ws.Cells.Item(2,7).Value = 2

# New row 3 - REYES PEREZ ARANTZA (2BLCM)
$ws.Cells.Item(3,1).Value = 20330051920310
$ws.Cells.Item(3,2).Value = "REYES"
$ws.Cells.Item(3,3).Value = "PEREZ"
$ws.Cells.Item(3,4).Value = "ARANTZA"
$ws.Cells.Item(3,5).Value = "TOMA MUESTRAS BIOLÓGICAS"
$ws.Cells.Item(3,6).Value = "2BLCM"
$ws.Cells.Item(3,7).Value = 2

# Row 4 - previously row 2 (MARIANO ANTONIO JAQUELINE, 2ALCM); same
# student/NC but Reprobadas corrected from 2 to 1
$ws.Cells.Item(4,1).Value = 20330051920237
$ws.Cells.Item(4,2).Value = "MARIANO"
$ws.Cells.Item(4,3).Value = "ANTONIO"
$ws.Cells.Item(4,4).Value = "JAQUELINE"
$ws.Cells.Item(4,5).Value = "TOMA MUESTRAS BIOLÓGICAS"
$ws.Cells.Item(4,6).Value = "2ALCM"
$ws.Cells.Item(4,7).Value = 1
